$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (shared strings content change):
#   B1: "N"  -> "N, unit"
#   C1: "R"  -> "R_ev, Om"
#   D1: "R"  -> "R, Om"
$ws.Range("B1").Value = "N, unit"
$ws.Range("C1").Value = "R_ev, Om"
$ws.Range("D1").Value = "R, Om"

# Move the selection from E4 to B8
$ws.Range("B8").Select()
